$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data table (rows 2-9, columns A-T).
# Rows 2-5: "ECs" sending-cluster variant of the Cd5l-Cd5 edge (previously "M2").
# Rows 6-9: new "M2" sending-cluster rows for the same edge (appended after row 5).
$data = @(
    @("ECs", "Cd5l", "Cd5", "ECs",  1, 0.3333333333333333, 0.8129133333333334, 2.43874,  0.3118762449524117, 0.3118762449524117, 3, 1, 1.992340666666667,  5.977022,  0.5190478663364179, 0.5190478663364179, 1.619600292475555,  14.57640263228,      0.1618786995035633,   0.1618786995035633),
    @("ECs", "Cd5l", "Cd5", "FAPs", 1, 0.3333333333333333, 0.8129133333333334, 2.43874,  0.3118762449524117, 0.3118762449524117, 3, 1, 0.519652,            1.558956,  0.1353805934648319, 0.1353805934648319, 0.4224320394933334, 3.80188835544,       0.04222199112924079,  0.0422219911292408),
    @("ECs", "Cd5l", "Cd5", "M2",   1, 0.3333333333333333, 0.8129133333333334, 2.43874,  0.3118762449524117, 0.3118762449524117, 3, 1, 0.8421986666666667,  2.526596,  0.2194109814041387, 0.2194109814041387, 0.6846345254488889, 6.16171072904,       0.06842907298164623,  0.06842907298164623),
    @("ECs", "Cd5l", "Cd5", "sCs",  1, 0.3333333333333333, 0.8129133333333334, 2.43874,  0.3118762449524117, 0.3118762449524117, 3, 1, 0.4842613333333333,  1.452784,  0.1261605587946115, 0.1261605587946115, 0.3936624946844444, 3.54296245216,       0.0393464813379614,   0.0393464813379614),
    @("M2",  "Cd5l", "Cd5", "ECs",  3, 1,                  1.793612,            5.380836, 0.6881237550475882, 0.6881237550475883, 3, 1, 1.992340666666667,  5.977022,  0.5190478663364179, 0.5190478663364179, 3.573486127821333,  32.161375150392,     0.3571691668328545,   0.3571691668328545),
    @("M2",  "Cd5l", "Cd5", "FAPs", 3, 1,                  1.793612,            5.380836, 0.6881237550475882, 0.6881237550475883, 3, 1, 0.519652,            1.558956,  0.1353805934648319, 0.1353805934648319, 0.9320540630240001, 8.388486567216,      0.09315860233559112,  0.09315860233559115),
    @("M2",  "Cd5l", "Cd5", "M2",   3, 1,                  1.793612,            5.380836, 0.6881237550475882, 0.6881237550475883, 3, 1, 0.8421986666666667,  2.526596,  0.2194109814041387, 0.2194109814041387, 1.510577634917333,  13.595198714256,     0.1509819084224925,   0.1509819084224925),
    @("M2",  "Cd5l", "Cd5", "sCs",  3, 1,                  1.793612,            5.380836, 0.6881237550475882, 0.6881237550475883, 3, 1, 0.4842613333333333,  1.452784,  0.1261605587946115, 0.1261605587946115, 0.8685769386026666, 7.817192447424,      0.08681407745665008,  0.0868140774566501)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $row[$j]
    }
}
